$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D7: 283 -> 373 (this also drives the L2 AVERAGE(D2:D19) recalculation)
$ws.Range("D7").Value = 373

# Row 8 height: 15 -> 13.8
$ws.Rows.Item(8).RowHeight = 13.8

# New cell H8: date 45860 (22/07/2025), same date style as H2:H4
$ws.Range("H8").Value = 45860
$ws.Range("H8").NumberFormat = "dd/mm/yyyy"

# G11: 568 -> 752
$ws.Range("G11").Value = 752

# New cell E18: 0
$ws.Range("E18").Value = 0

# Update the active selection from F7 to D6
$ws.Range("D6").Select()
